# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Updates column G ("K") values for rows 2-40 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newKValues = @{
    2  = 3
    3  = 8
    4  = 4
    5  = 3
    6  = 4
    7  = 4
    8  = 6
    9  = 7
    10 = 8
    11 = 4
    12 = 6
    13 = 4
    14 = 4
    15 = 4
    16 = 6
    17 = 6
    18 = 8
    19 = 5
    20 = 5
    21 = 6
    22 = 6
    23 = 8
    24 = 6
    25 = 5
    26 = 6
    27 = 6
    28 = 4
    29 = 4
    30 = 5
    31 = 4
    32 = 3
    33 = 5
    34 = 5
    35 = 3
    36 = 8
    37 = 1
    38 = 6
    39 = 2
    40 = 2
}

foreach ($row in $newKValues.Keys) {
    $ws.Range("G$row").Value = $newKValues[$row]
}
